# Apply the edit: insert a new "Paper Name" header column before the
# "Paper Language" column on Sheet1, with its own highlighted/bordered style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at F (shifts Paper Language, Number of Questions,
# Exam Date one column to the right: F->G, G->H, H->I)
$ws.Columns("F:F").Insert()

# Set the new header cell content
$ws.Range("F1").Value = "Paper Name"

# Style the new header cell: centered, vertically centered, wrap text,
# bold red Arial Narrow font, light blue fill, thin black border around.
$ws.Range("F1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F1").VerticalAlignment = -4108     # xlCenter
$ws.Range("F1").WrapText = $true

$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").Font.Size = 12
$ws.Range("F1").Font.Name = "Arial Narrow"
$ws.Range("F1").Font.Color = 255

$ws.Range("F1").Interior.Pattern = 1  # xlSolid
$ws.Range("F1").Interior.ThemeColor = 5   # msoThemeColorAccent1
$ws.Range("F1").Interior.TintAndShade = 0.79998168889431442

$ws.Range("F1").Borders.LineStyle = 1
$ws.Range("F1").Borders.Weight = 2

$ws.Range("A1").Select()
